# Insert a new data row at row 32 (shifting the existing rows 32..103 down
# to 33..104) and populate it with a new "Cilantro" price record for the
# "Macroferia Regional de Talca" market.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push existing row 32 (and everything below it) down by one row.
$ws.Rows(32).Insert()

# Fill in the newly inserted row with the new observation.
$ws.Range("A32").Value = 5
$ws.Range("B32").Value = "Macroferia Regional de Talca"
$ws.Range("C32").Value = "Maule"
$ws.Range("D32").Value = 45133
$ws.Range("E32").Value = 7
$ws.Range("F32").Value = 100112040
$ws.Range("G32").Value = "Cilantro"
$ws.Range("H32").Value = "Sin especificar"
$ws.Range("I32").Value = "Primera"
$ws.Range("J32").Value = 150
$ws.Range("K32").Value = 9000
$ws.Range("L32").Value = 9000
$ws.Range("M32").Value = 9000
$ws.Range("N32").Value = "$/caja 36 atados"
$ws.Range("O32").Value = "Región Metropolitana"
$ws.Range("P32").Value = 250
$ws.Range("Q32").Value = 36
$ws.Range("R32").Value = "Hortaliza"
